# Rename the "SWC" sheet to "SWC_Composition".
$wb = $excel.ActiveWorkbook
$swc = $wb.Worksheets.Item("SWC")
$swc.Name = "SWC_Composition"

# Update the selection left on the (now renamed) SWC_Composition sheet
# while it is briefly made active, matching the recorded view state
# (selection moved from E30 to E32, no longer the active tab).
$swc.Activate()
$swc.Range("E32").Select()

# Finish with "DataTypes" as the active sheet/tab, with its selection
# moved from H23 to I21 (matches tabSelected on sheet1 + activeTab=0
# on the workbook).
$dataTypes = $wb.Worksheets.Item("DataTypes")
$dataTypes.Activate()
$dataTypes.Range("I21").Select()
